# feat: implement 13 PRSB alignment profiles and terminology
#
# Regenerates this ValueSet page for the updated "goal-evaluation-valueset"
# terminology: bumps the generation timestamp on the Metadata sheet and
# swaps the SNOMED CT concept list on the "Include #0" sheet from the old
# 3-code list to the new 5-code list (pushing the trailing blank row and
# the "System URI" row down to make room).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date -----------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-01T13:37:23+00:00"

# --- Include #0 sheet: refresh the Concept/Description rows ---------------
$ws = $wb.Worksheets.Item("Include #0")

# SNOMED CT concept codes are textual identifiers, not numbers - format the
# column as text up front so the new codes aren't auto-coerced to numerics.
$ws.Range("A2").NumberFormat = "@"

# There were 3 concept rows (rows 2-4), then a blank spacer row (row 5) and
# a "System URI" row (row 6). The new concept list has 5 rows, so insert two
# extra rows before the spacer to make room, carrying the spacer/System URI
# rows down to rows 7-8.
$ws.Rows("5:6").Insert(-4121)

# Insert() only loosely approximates the row-above formatting for the new
# rows; make sure rows 3-6 exactly match row 2's (text column + border/fill).
$ws.Range("A2:B2").Copy()
$ws.Range("A3:B6").PasteSpecial(-4122)

$ws.Range("A2").Value = "385652002"
$ws.Range("B2").Value = "Objective achieved"

$ws.Range("A3").Value = "385651009"
$ws.Range("B3").Value = "Objective not achieved"

$ws.Range("A4").Value = "255609007"
$ws.Range("B4").Value = "Partial achievement"

$ws.Range("A5").Value = "723510000"
$ws.Range("B5").Value = "Sustained improvement"

$ws.Range("A6").Value = "260388008"
$ws.Range("B6").Value = "Worsening"
